# The "Recorded By" column (G) holds a comma-separated list of the users who
# recorded/updated each attendance session (e.g. "System, dnasr281@gmail.com").
# This edit reverses the order of the comma-separated entries for every row
# whose list currently starts with "System" (e.g. "System, X" -> "X, System",
# "System, system, backup@backdoor.com" -> "backup@backdoor.com, system, System").
# Rows whose "Recorded By" value is a single entry, or that do not start with
# "System", are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null) {
        $valStr = $val.ToString()

        if ($valStr.Contains(",")) {
            $parts = $valStr.Split(",")
            $firstEntry = $parts[0].Trim()

            if ($firstEntry -eq "System") {
                $count = $parts.Count
                $reversedParts = @()
                for ($i = $count - 1; $i -ge 0; $i--) {
                    $reversedParts += $parts[$i].Trim()
                }
                $newVal = [string]::Join(", ", $reversedParts)
                $cell.Value = $newVal
            }
        }
    }
}
